# Add working set of sequences
# Replace the pool of German verb stimuli in column B (rows 2-193) with the
# new working set, in order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "hauen", "wundern", "liefern", "gelten", "münzen", "dürfen", "bauen", "schlagen",
    "bluten", "geben", "töten", "zielen", "stammen", "sichern", "treffen", "schwören",
    "sprechen", "warnen", "wehen", "quälen", "schrecken", "landen", "altern", "tropfen",
    "helfen", "trauen", "deuten", "dringen", "wüten", "betteln", "zögern", "schwächen",
    "heißen", "wirken", "kümmern", "scheinen", "leiden", "enden", "boxen", "fügen",
    "schreiten", "flehen", "spüren", "leeren", "saufen", "irren", "biegen", "mögen",
    "wärmen", "süßen", "fangen", "planen", "lügen", "lassen", "suchen", "siegen",
    "retten", "folgen", "ändern", "räumen", "hören", "kaufen", "streichen", "stürzen",
    "locken", "wehtun", "heulen", "klettern", "meistern", "schenken", "malen", "werfen",
    "schlucken", "äußern", "pflanzen", "handeln", "knarren", "jagen", "brauchen", "achten",
    "wachsen", "buchen", "foltern", "küssen", "feiern", "heben", "treten", "platzen",
    "reiten", "reizen", "bitten", "passen", "reisen", "fließen", "albern", "treiben",
    "schmecken", "pfeifen", "starten", "sprengen", "greifen", "spielen", "bieten", "schulden",
    "arten", "ärgern", "nerven", "sperren", "eignen", "sorgen", "fällen", "rasen",
    "bergen", "wählen", "grüßen", "baden", "mauern", "hacken", "beißen", "sterben",
    "kichern", "weichen", "decken", "dienen", "graben", "fragen", "fahren", "gründen",
    "freuen", "heilen", "bremsen", "bellen", "schließen", "spenden", "tollen", "seufzen",
    "drehen", "atmen", "parken", "formen", "stillen", "liegen", "klingen", "stecken",
    "machen", "führen", "hupen", "spinnen", "scheiden", "werden", "weigern", "wollen",
    "trotzen", "laufen", "kosten", "streifen", "fischen", "drucken", "loben", "zählen",
    "schwingen", "wenden", "wetten", "zünden", "segnen", "beten", "fallen", "posten",
    "leihen", "backen", "filmen", "stehlen", "lockern", "knien", "zeigen", "kehren",
    "prügeln", "trennen", "runden", "erben", "füllen", "schützen", "sinken", "scheitern",
    "lesen", "rufen", "ehren", "haben", "stärken", "jubeln", "ziehen", "flüchten"
)

$startRow = 2
for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $words[$i]
}
